$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: add new columns G (Ny leverandør) and H (TCV_range),
#     and rename the existing F header from "TCV_range" to "Årsag" ---

# Give the new header cells the same look as the existing bold/bordered
# header cells (copy format from F1, which already carries the header style).
$ws.Range("F1").Copy() | Out-Null
$ws.Range("G1:H1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("F1").Value = "Årsag"
$ws.Range("G1").Value = "Ny leverandør"
$ws.Range("H1").Value = "TCV_range"

# --- Data rows: move the old "TCV_range" value (160000-180000) from column F
#     to the new column H, and populate column F with the reason/"Årsag" ---

$ws.Range("H2").Value = "160000-180000"
$ws.Range("H3").Value = "160000-180000"
$ws.Range("H4").Value = "160000-180000"
$ws.Range("H5").Value = "160000-180000"
$ws.Range("H6").Value = "160000-180000"

$ws.Range("F2").Value = "Utilfredshed (Service - uddyb i bemærkninger)"
$ws.Range("F3").Value = "Strategisk beslutning"
$ws.Range("F4").Value = "Ikke oplyst"
$ws.Range("F5").Value = "Ikke oplyst"
$ws.Range("F6").Value = "Utilfredshed (Service - uddyb i bemærkninger)"
